$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{
        A = "Golang Developer-67024"
        B = "https://www.dice.com/job-detail/18cd339e-d2db-43d3-aa6d-1e353907b413?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
        C = "San Jose, California"
        D = "Contract"
        E = "Depends on Experience"
        F = "Keypixel Software Solutions"
    },
    @{
        A = "Golang Developer"
        B = "https://www.dice.com/job-detail/c598b4f5-0a1a-464e-b177-b2e8c419645a?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
        C = "Virginia"
        D = "Contract"
        E = "USD85 - USD110"
        F = "Dexian DISYS"
    },
    @{
        A = "Golang Developer – AWS & Microservices"
        B = "https://www.dice.com/job-detail/5b12ce56-e8f0-4287-87ca-6c1f607ee0a9?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
        C = "McLean, Virginia"
        D = "Contract"
        E = "55 - 60"
        F = "Rapsys Technologies"
    },
    @{
        A = "Software Development Engineer (GoLang)- Onsite"
        B = "https://www.dice.com/job-detail/9d97cc6a-0bb3-4668-8749-385f92ff283b?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
        C = "Chantilly, Virginia"
        D = "Contract"
        E = "Depends on Experience"
        F = "Stellar Professionals LLC"
    }
)

$startRow = 181
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
}
